$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTechFramework")

# Remove the "Coal Steam" row (old row 2): rows below shift up and
# formulas referencing F10/F11 auto-adjust to F9/F10 etc.
$ws.Rows.Item(2).Delete()

# Restore the final selection state recorded in the saved file.
$ws.Range("C8").Select()
